$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.525996016902489
$ws.Range("D2").Value = 4.16254988427642
$ws.Range("E2").Value = 11.37259764713656
$ws.Range("F2").Value = 19.24634769706535
$ws.Range("G2").Value = 20.04906349981956
$ws.Range("H2").Value = 12.14333230258115
$ws.Range("I2").Value = 16.47740128438427
$ws.Range("K2").Value = 13.00578605792526
$ws.Range("M2").Value = 14.78986138790676
$ws.Range("O2").Value = 17.30873502555927
$ws.Range("C3").Value = 3.359023480698185
$ws.Range("D3").Value = 4.069819882986558
$ws.Range("E3").Value = 11.3489573283096
$ws.Range("F3").Value = 19.32001459102019
$ws.Range("G3").Value = 20.16986236300911
$ws.Range("H3").Value = 12.20724497716377
$ws.Range("I3").Value = 16.55918996104725
$ws.Range("K3").Value = 12.3265762931016
$ws.Range("M3").Value = 14.41950991951688
$ws.Range("O3").Value = 17.41864455020784
$ws.Range("C4").Value = 3.251436890617226
$ws.Range("D4").Value = 4.011309235506118
$ws.Range("E4").Value = 11.34043073973271
$ws.Range("F4").Value = 19.37282834446628
$ws.Range("G4").Value = 20.25555007798444
$ws.Range("H4").Value = 12.24910654776711
$ws.Range("I4").Value = 16.6147715577058
$ws.Range("K4").Value = 11.88808443033535
$ws.Range("M4").Value = 14.18925635101785
$ws.Range("O4").Value = 17.49160899335379
$ws.Range("C5").Value = 3.206352324261267
$ws.Range("D5").Value = 3.987091323915922
$ws.Range("E5").Value = 11.3384586422425
$ws.Range("F5").Value = 19.39624531015028
$ws.Range("G5").Value = 20.29333359571513
$ws.Range("H5").Value = 12.26682363331617
$ws.Range("I5").Value = 16.6387637036264
$ws.Range("K5").Value = 11.7041247233049
$ws.Range("M5").Value = 14.09484201819463
$ws.Range("O5").Value = 17.52271505868261
$ws.Range("C6").Value = 3.198792070777751
$ws.Range("D6").Value = 3.983047999873351
$ws.Range("E6").Value = 11.33822177922922
$ws.Range("F6").Value = 19.40024777569609
$ws.Range("G6").Value = 20.29977963915553
$ws.Range("H6").Value = 12.26980528189319
$ws.Range("I6").Value = 16.64282845748711
$ws.Range("K6").Value = 11.67326413552232
$ws.Range("M6").Value = 14.07913335436005
$ws.Range("O6").Value = 17.52796293321363
$ws.Range("C7").Value = 3.250833848272992
$ws.Range("D7").Value = 4.010984110746596
$ws.Range("E7").Value = 11.34039806590137
$ws.Range("F7").Value = 19.37313649722569
$ws.Range("G7").Value = 20.25604808080061
$ws.Range("H7").Value = 12.24934282224212
$ws.Range("I7").Value = 16.61508969822055
$ws.Range("K7").Value = 11.88562464412479
$ws.Range("M7").Value = 14.18798522953688
$ws.Range("O7").Value = 17.49202295071762
$ws.Range("C8").Value = 3.469492276688216
$ws.Range("D8").Value = 4.130915760143463
$ws.Range("E8").Value = 11.36320230653562
$ws.Range("F8").Value = 19.27016762611894
$ws.Range("G8").Value = 20.08830642570111
$ws.Range("H8").Value = 12.1648255915589
$ws.Range("I8").Value = 16.50448512798771
$ws.Range("K8").Value = 12.77611449378603
$ws.Range("M8").Value = 14.66283542269985
$ws.Range("O8").Value = 17.3454912815921
$ws.Range("C9").Value = 3.856972569890714
$ws.Range("D9").Value = 4.352734292586979
$ws.Range("E9").Value = 11.45545959636342
$ws.Range("F9").Value = 19.12889126563325
$ws.Range("G9").Value = 19.85206458384134
$ws.Range("H9").Value = 12.01988547024816
$ws.Range("I9").Value = 16.33042240944805
$ws.Range("K9").Value = 14.34781733359237
$ws.Range("M9").Value = 15.56580992271994
$ws.Range("O9").Value = 17.10186451071581
$ws.Range("C10").Value = 4.115277364655967
$ws.Range("D10").Value = 4.506473066620459
$ws.Range("E10").Value = 11.55206601624899
$ws.Range("F10").Value = 19.06267780476875
$ws.Range("G10").Value = 19.73670956819512
$ws.Range("H10").Value = 11.92610179872917
$ws.Range("I10").Value = 16.2290230898082
$ws.Range("K10").Value = 15.39167745792539
$ws.Range("M10").Value = 16.20514208943863
$ws.Range("O10").Value = 16.94986133836161
$ws.Range("C11").Value = 4.226879677429754
$ws.Range("D11").Value = 4.574205565384279
$ws.Range("E11").Value = 11.60218436813628
$ws.Range("F11").Value = 19.04082886604363
$ws.Range("G11").Value = 19.69720033347084
$ws.Range("H11").Value = 11.8862017806187
$ws.Range("I11").Value = 16.1887212249527
$ws.Range("K11").Value = 15.84184227580525
$ws.Range("M11").Value = 16.48946479369715
$ws.Range("O11").Value = 16.88664483862645
$ws.Range("C12").Value = 4.268280570108254
$ws.Range("D12").Value = 4.599521930949057
$ws.Range("E12").Value = 11.62203887818198
$ws.Range("F12").Value = 19.03375247664969
$ws.Range("G12").Value = 19.68412786532173
$ws.Range("H12").Value = 11.87149060064577
$ws.Range("I12").Value = 16.17430363030256
$ws.Range("K12").Value = 16.00871861720451
$ws.Range("M12").Value = 16.59609694435946
$ws.Range("O12").Value = 16.8635654028319
$ws.Range("C13").Value = 4.259402580009551
$ws.Range("D13").Value = 4.594084618681218
$ws.Range("E13").Value = 11.61772411631669
$ws.Range("F13").Value = 19.03522313841462
$ws.Range("G13").Value = 19.68685888630307
$ws.Range("H13").Value = 11.87464119769208
$ws.Range("I13").Value = 16.17737109077388
$ws.Range("K13").Value = 15.97293908230298
$ws.Range("M13").Value = 16.57317929635891
$ws.Range("O13").Value = 16.86849765188357
$ws.Range("C14").Value = 4.230303058263712
$ws.Range("D14").Value = 4.576295105453521
$ws.Range("E14").Value = 11.60380033454349
$ws.Range("F14").Value = 19.04022264975267
$ws.Range("G14").Value = 19.69608686761489
$ws.Range("H14").Value = 11.88498349978159
$ws.Range("I14").Value = 16.18751813415259
$ws.Range("K14").Value = 15.85564346698507
$ws.Range("M14").Value = 16.49825875459101
$ws.Range("O14").Value = 16.8847288210154
$ws.Range("C15").Value = 4.212366391460066
$ws.Range("D15").Value = 4.565354779564028
$ws.Range("E15").Value = 11.59538527875581
$ws.Range("F15").Value = 19.04344113741955
$ws.Range("G15").Value = 19.7019859360391
$ws.Range("H15").Value = 11.89137032682016
$ws.Range("I15").Value = 16.19384356490756
$ws.Range("K15").Value = 15.78332777896312
$ws.Range("M15").Value = 16.45223023991054
$ws.Range("O15").Value = 16.89478297756261
$ws.Range("C16").Value = 4.107864066238025
$ws.Range("D16").Value = 4.502000881532207
$ws.Range("E16").Value = 11.54891391049027
$ws.Range("F16").Value = 19.06427282984021
$ws.Range("G16").Value = 19.73955463035393
$ws.Range("H16").Value = 11.9287650199281
$ws.Range("I16").Value = 16.23177466353805
$ws.Range("K16").Value = 15.36175756603675
$ws.Range("M16").Value = 16.18642129952337
$ws.Range("O16").Value = 16.95411257525761
$ws.Range("C17").Value = 4.042234486574701
$ws.Range("D17").Value = 4.462559358932499
$ws.Range("E17").Value = 11.52197841079821
$ws.Range("F17").Value = 19.07917693437948
$ws.Range("G17").Value = 19.76594075382286
$ws.Range("H17").Value = 11.95241360770577
$ws.Range("I17").Value = 16.2565406104676
$ws.Range("K17").Value = 15.09678159644054
$ws.Range("M17").Value = 16.02161559111843
$ws.Range("O17").Value = 16.99203290387435
$ws.Range("C18").Value = 4.00393111585186
$ws.Range("D18").Value = 4.439667029490402
$ws.Range("E18").Value = 11.50706758828217
$ws.Range("F18").Value = 19.08852756598222
$ws.Range("G18").Value = 19.7823364864561
$ws.Range("H18").Value = 11.96627563118025
$ws.Range("I18").Value = 16.27133326602815
$ws.Range("K18").Value = 14.94205201084955
$ws.Range("M18").Value = 15.92621604035087
$ws.Range("O18").Value = 17.01440149058712
$ws.Range("C19").Value = 3.990867335630579
$ws.Range("D19").Value = 4.431881088726282
$ws.Range("E19").Value = 11.50211927357829
$ws.Range("F19").Value = 19.09182692743839
$ws.Range("G19").Value = 19.78809641193735
$ws.Range("H19").Value = 11.9710137016726
$ws.Range("I19").Value = 16.27643572167587
$ws.Range("K19").Value = 14.88926575020954
$ws.Range("M19").Value = 15.89381415103218
$ws.Range("O19").Value = 17.02207071852369
$ws.Range("C20").Value = 4.049278398088758
$ws.Range("D20").Value = 4.466779476058216
$ws.Range("E20").Value = 11.52478560700695
$ws.Range("F20").Value = 19.07750976769799
$ws.Range("G20").Value = 19.76300554752647
$ws.Range("H20").Value = 11.94986926220965
$ws.Range("I20").Value = 16.25384747805525
$ws.Range("K20").Value = 15.12522942621447
$ws.Range("M20").Value = 16.03922306858586
$ws.Range("O20").Value = 16.98793844139902
$ws.Range("C21").Value = 4.238873728663521
$ws.Range("D21").Value = 4.58152945705439
$ws.Range("E21").Value = 11.60786642391585
$ws.Range("F21").Value = 19.03872161957773
$ws.Range("G21").Value = 19.69332494116821
$ws.Range("H21").Value = 11.88193490410315
$ws.Range("I21").Value = 16.1845147515556
$ws.Range("K21").Value = 15.89019377627561
$ws.Range("M21").Value = 16.5202935564249
$ws.Range("O21").Value = 16.87993796069827
$ws.Range("C22").Value = 4.357765826842588
$ws.Range("D22").Value = 4.654581046957725
$ws.Range("E22").Value = 11.66726209349898
$ws.Range("F22").Value = 19.02035316754856
$ws.Range("G22").Value = 19.65880380575856
$ws.Range("H22").Value = 11.83985666402308
$ws.Range("I22").Value = 16.14412310191613
$ws.Range("K22").Value = 16.3691968097448
$ws.Range("M22").Value = 16.82863021842789
$ws.Range("O22").Value = 16.81436498429354
$ws.Range("C23").Value = 4.294773453159526
$ws.Range("D23").Value = 4.61577473172752
$ws.Range("E23").Value = 11.63509957924062
$ws.Range("F23").Value = 19.0295155977066
$ws.Range("G23").Value = 19.67621268708906
$ws.Range("H23").Value = 11.86210197704415
$ws.Range("I23").Value = 16.16522860026602
$ws.Range("K23").Value = 16.11547142701681
$ws.Range("M23").Value = 16.66465092189916
$ws.Range("O23").Value = 16.84890180397199
$ws.Range("C24").Value = 4.046095626599747
$ws.Range("D24").Value = 4.464872235423191
$ws.Range("E24").Value = 11.52351468243084
$ws.Range("F24").Value = 19.07826105828009
$ws.Range("G24").Value = 19.76432873725185
$ws.Range("H24").Value = 11.95101873185791
$ws.Range("I24").Value = 16.25506331681123
$ws.Range("K24").Value = 15.11237560784203
$ws.Range("M24").Value = 16.03126474680554
$ws.Range("O24").Value = 16.9897877800473
$ws.Range("C25").Value = 3.756709446396635
$ws.Range("D25").Value = 4.294274398983971
$ws.Range("E25").Value = 11.42541772307648
$ws.Range("F25").Value = 19.16055131490173
$ws.Range("G25").Value = 19.90586849478484
$ws.Range("H25").Value = 12.05686662789099
$ws.Range("I25").Value = 16.37288635150249
$ws.Range("K25").Value = 13.94183987351289
$ws.Range("M25").Value = 15.32529306826524
$ws.Range("O25").Value = 17.16305689228697
